# Updated cryptos list on Sun Jul 28 04:42:08 UTC 2024 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) columns of the
# cryptos sheet with newly scraped figures. Only the rows whose price and/or
# percentage actually moved are touched; everything else (Coin name, Link,
# rank index) is left exactly as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new Price (column D) text. Left out of this table entirely when the
# price column did not change for that row.
$priceUpdates = [ordered]@{
    2  = "67.735.63"
    3  = "3.235.48"
    5  = "579.53"
    6  = "183.67"
    11 = "0.413"
    12 = "3.798.70"
    15 = "67.746.22"
    17 = "3.223.14"
    19 = "13.42"
    20 = "395.66"
    21 = "7.54"
    22 = "1.00"
    23 = "71.05"
    26 = "0.185"
    27 = "9.51"
    30 = "5.56"
    31 = "22.60"
    35 = "161.75"
    37 = "1.88"
    38 = "0.807"
    39 = "26.32"
    42 = "41.09"
    45 = "2.604.63"
    46 = "24.81"
    47 = "334.92"
    51 = "30.83"
}

# row -> new Volume(1h) (column E) text (the two leading/trailing spaces are
# part of the original cell content).
$volumeUpdates = [ordered]@{
    2  = "  -0.11%  "
    3  = "  -0.52%  "
    4  = "  -0.01%  "
    5  = "  -0.70%  "
    6  = "  +0.23%  "
    7  = "  +0.00%  "
    8  = "  +0.22%  "
    9  = "  -3.83%  "
    10 = "  -1.35%  "
    11 = "  -0.33%  "
    12 = "  -0.48%  "
    13 = "  +0.03%  "
    14 = "  -3.52%  "
    15 = "  -0.09%  "
    16 = "  -1.60%  "
    17 = "  -0.62%  "
    18 = "  -1.07%  "
    19 = "  -1.08%  "
    20 = "  +3.80%  "
    21 = "  -1.34%  "
    22 = "  -0.06%  "
    23 = "  -0.23%  "
    24 = "  -0.02%  "
    25 = "  -1.56%  "
    26 = "  +2.57%  "
    27 = "  -3.04%  "
    28 = "  +0.04%  "
    29 = "  -1.39%  "
    30 = "  -1.87%  "
    31 = "  -1.23%  "
    32 = "  -2.35%  "
    33 = "  -1.20%  "
    35 = "  +0.25%  "
    36 = "  -4.78%  "
    37 = "  +1.64%  "
    38 = "  -3.44%  "
    39 = "  -0.95%  "
    40 = "  -1.29%  "
    41 = "  -3.27%  "
    42 = "  -0.44%  "
    43 = "  -4.79%  "
    44 = "  -0.81%  "
    45 = "  -1.09%  "
    46 = "  -2.51%  "
    47 = "  -3.55%  "
    48 = "  -2.19%  "
    49 = "  +0.94%  "
    50 = "  -2.10%  "
    51 = "  +0.52%  "
}

# Plain numeric-looking text ("71.05", "1.00", ...) would otherwise be
# auto-converted to a real number by Excel on assignment. Prefixing with an
# apostrophe forces it to stay text, matching the source data (which stores
# every Price cell as text, including ones such as "67.735.63" that contain
# more than one '.' and so are never at risk of numeric auto-conversion).
function Set-TextValue($range, [string]$text) {
    if ($text -match '^[+-]?\d+(\.\d+)?$') {
        $range.Value = "'" + $text
    } else {
        $range.Value = $text
    }
}

foreach ($row in $priceUpdates.Keys) {
    Set-TextValue $ws.Range("D$row") $priceUpdates[$row]
}

foreach ($row in $volumeUpdates.Keys) {
    Set-TextValue $ws.Range("E$row") $volumeUpdates[$row]
}
